$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range("D2").Value = "65.683.86"
$ws.Range("E2").Value = "  -5.33%  "

$ws.Range("D3").Value = "3.472.75"
$ws.Range("E3").Value = "  -7.13%  "

Set-TextValue $ws.Range("D4") "0.999"
$ws.Range("E4").Value = "  -0.05%  "

Set-TextValue $ws.Range("D5") "559.82"
$ws.Range("E5").Value = "  -8.72%  "

Set-TextValue $ws.Range("D6") "181.88"
$ws.Range("E6").Value = "  -5.86%  "

$ws.Range("D7").Value = "3.467.04"
$ws.Range("E7").Value = "  -7.13%  "

Set-TextValue $ws.Range("D8") "0.600"
$ws.Range("E8").Value = "  -6.13%  "

Set-TextValue $ws.Range("D9") "1.00"
$ws.Range("E9").Value = "  +0.04%  "

Set-TextValue $ws.Range("D10") "0.648"
$ws.Range("E10").Value = "  -11.61%  "

Set-TextValue $ws.Range("D11") "0.141"
$ws.Range("E11").Value = "  -13.49%  "

Set-TextValue $ws.Range("D12") "51.35"
$ws.Range("E12").Value = "  -14.99%  "

Set-TextValue $ws.Range("D13") "0.0000251"
$ws.Range("E13").Value = "  -14.47%  "

Set-TextValue $ws.Range("D14") "9.48"
$ws.Range("E14").Value = "  -11.13%  "

$ws.Range("D15").Value = "4.020.54"
$ws.Range("E15").Value = "  -7.19%  "

Set-TextValue $ws.Range("D16") "0.125"
$ws.Range("E16").Value = "  -1.78%  "

$ws.Range("D17").Value = "3.464.82"
$ws.Range("E17").Value = "  -7.16%  "

$ws.Range("D18").Value = "65.377.82"
$ws.Range("E18").Value = "  -5.53%  "

Set-TextValue $ws.Range("D19") "17.67"
$ws.Range("E19").Value = "  -9.73%  "

Set-TextValue $ws.Range("D20") "11.66"
$ws.Range("E20").Value = "  -10.27%  "

Set-TextValue $ws.Range("D21") "1.03"
$ws.Range("E21").Value = "  -10.97%  "

Set-TextValue $ws.Range("D22") "377.64"
$ws.Range("E22").Value = "  -8.97%  "

Set-TextValue $ws.Range("D23") "4.08"
$ws.Range("E23").Value = "  -11.18%  "

Set-TextValue $ws.Range("D24") "82.87"
$ws.Range("E24").Value = "  -8.05%  "

Set-TextValue $ws.Range("D25") "10.73"
$ws.Range("E25").Value = "  -3.57%  "

$ws.Range("B26").Value = "ImmutableX"
$ws.Range("C26").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws.Range("D26") "2.81"
$ws.Range("E26").Value = "  -9.17%  "

$ws.Range("B27").Value = "LEO"
$ws.Range("C27").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue $ws.Range("D27") "5.99"
$ws.Range("E27").Value = "  -1.16%  "

Set-TextValue $ws.Range("D28") "11.84"
$ws.Range("E28").Value = "  -8.81%  "

Set-TextValue $ws.Range("D29") "3.44"
$ws.Range("E29").Value = "  -10.21%  "

Set-TextValue $ws.Range("D30") "8.59"
$ws.Range("E30").Value = "  -12.22%  "

Set-TextValue $ws.Range("D31") "30.23"
$ws.Range("E31").Value = "  -9.22%  "

Set-TextValue $ws.Range("D32") "7.21"
$ws.Range("E32").Value = "  -8.12%  "

Set-TextValue $ws.Range("D33") "608.64"
$ws.Range("E33").Value = "  -4.18%  "

Set-TextValue $ws.Range("D34") "11.84"
$ws.Range("E34").Value = "  -7.82%  "

Set-TextValue $ws.Range("D35") "62.52"
$ws.Range("E35").Value = "  -6.74%  "

$ws.Range("E36").Value = "  -11.83%  "

Set-TextValue $ws.Range("D37") "40.63"
$ws.Range("E37").Value = "  -11.49%  "

$ws.Range("E38").Value = "  +0.39%  "

Set-TextValue $ws.Range("D39") "0.393"
$ws.Range("E39").Value = "  -6.15%  "

Set-TextValue $ws.Range("D40") "0.996"
$ws.Range("E40").Value = "  -0.39%  "

$ws.Range("D41").Value = "0.0₃0712"
$ws.Range("E41").Value = "  -15.83%  "

Set-TextValue $ws.Range("D42") "0.129"
$ws.Range("E42").Value = "  -9.30%  "

$ws.Range("D43").Value = "2.902.75"
$ws.Range("E43").Value = "  +0.05%  "

Set-TextValue $ws.Range("D44") "2.73"
$ws.Range("E44").Value = "  -12.36%  "

$ws.Range("E45").Value = "  -8.88%  "

Set-TextValue $ws.Range("D46") "3.11"
$ws.Range("E46").Value = "  +0.66%  "

Set-TextValue $ws.Range("D47") "0.0394"
$ws.Range("E47").Value = "  -12.79%  "

$ws.Range("E48").Value = "  -10.05%  "

Set-TextValue $ws.Range("D49") "137.23"
$ws.Range("E49").Value = "  -4.84%  "

$ws.Range("B50").Value = "Stacks"
$ws.Range("C50").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue $ws.Range("D50") "2.66"
$ws.Range("E50").Value = "  -4.86%  "

Set-TextValue $ws.Range("D51") "8.11"
$ws.Range("E51").Value = "  -12.37%  "
